$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 144.23077
$ws.Range("I6").Value = 120.454544
$ws.Range("J6").Value = 275
$ws.Range("K6").Value = 361.363632
$ws.Range("L6").Value = 825
$ws.Range("M6").Value = -249.363632
$ws.Range("N6").Value = -1049
$ws.Range("H15").Value = 150.66
$ws.Range("I15").Value = 150.66
$ws.Range("K15").Value = 451.98
$ws.Range("M15").Value = -282.98
$ws.Range("H86").Value = 13415.75
$ws.Range("I86").Value = 900
$ws.Range("J86").Value = 25931.5
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 25931.5
$ws.Range("M86").Value = 223
$ws.Range("N86").Value = -28177.5
$ws.Range("H89").Value = 13415.75
$ws.Range("I89").Value = 900
$ws.Range("J89").Value = 25931.5
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 129657.5
$ws.Range("M89").Value = 1116
$ws.Range("N89").Value = -140889.5
$ws.Range("H118").Value = 600
$ws.Range("I118").Value = 600
$ws.Range("K118").Value = 1800
$ws.Range("M118").Value = -143
$ws.Range("H129").Value = 334319.16
$ws.Range("J129").Value = 435964.25
$ws.Range("L129").Value = 1307892.75
$ws.Range("N129").Value = -1317892.75
$ws.Range("H132").Value = 16668375
$ws.Range("I132").Value = 18183660
$ws.Range("J132").Value = 232.4
$ws.Range("K132").Value = 54550980
$ws.Range("L132").Value = 697.2
$ws.Range("M132").Value = -54548450
$ws.Range("N132").Value = -5757.2
$ws.Range("H138").Value = 4288.678
$ws.Range("I138").Value = 2969.3157
$ws.Range("J138").Value = 4657.3237
$ws.Range("K138").Value = 8907.947100000001
$ws.Range("L138").Value = 13971.9711
$ws.Range("M138").Value = -3767.947100000001
$ws.Range("N138").Value = -24251.9711
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2005
$ws.Range("I3").Value = 2005
$ws.Range("K3").Value = 2005
$ws.Range("M3").Value = -1890
$ws.Range("H32").Value = 18547.959
$ws.Range("I32").Value = 14920.647
$ws.Range("K32").Value = 14920.647
$ws.Range("M32").Value = -14633.647
$ws.Range("H61").Value = 13014294
$ws.Range("I61").Value = 17567792
$ws.Range("K61").Value = 17567792
$ws.Range("M61").Value = -17567580
$ws.Range("H88").Value = 168343.17
$ws.Range("I88").Value = 1766.3334
$ws.Range("J88").Value = 334920
$ws.Range("K88").Value = 1766.3334
$ws.Range("L88").Value = 334920
$ws.Range("M88").Value = -1360.3334
$ws.Range("N88").Value = -335732
$ws.Range("H91").Value = 168343.17
$ws.Range("I91").Value = 1766.3334
$ws.Range("J91").Value = 334920
$ws.Range("K91").Value = 1766.3334
$ws.Range("L91").Value = 334920
$ws.Range("M91").Value = -362.3334
$ws.Range("N91").Value = -337728
$ws.Range("H102").Value = 1020.4167
$ws.Range("I102").Value = 954.5
$ws.Range("J102").Value = 1350
$ws.Range("K102").Value = 954.5
$ws.Range("L102").Value = 1350
$ws.Range("M102").Value = 667.5
$ws.Range("N102").Value = -4594
$ws.Range("H136").Value = 13014294
$ws.Range("I136").Value = 17567792
$ws.Range("K136").Value = 52703376
$ws.Range("M136").Value = -52700826
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 400
$ws.Range("I11").Value = 400
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 400
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -260
$ws.Range("H42").Value = 119990
$ws.Range("J42").Value = 119990
$ws.Range("L42").Value = 119990
$ws.Range("N42").Value = -120646
$ws.Range("H86").Value = 2123.0454
$ws.Range("I86").Value = 1985.7142
$ws.Range("K86").Value = 1985.7142
$ws.Range("M86").Value = -862.7141999999999
$ws.Range("H89").Value = 2123.0454
$ws.Range("I89").Value = 1985.7142
$ws.Range("K89").Value = 9928.571
$ws.Range("M89").Value = -4312.571
$ws.Range("H94").Value = 846.2857
$ws.Range("I94").Value = 348.46155
$ws.Range("K94").Value = 348.46155
$ws.Range("M94").Value = 102.53845
$ws.Range("H99").Value = 1306.9
$ws.Range("I99").Value = 845
$ws.Range("J99").Value = 1999.75
$ws.Range("K99").Value = 845
$ws.Range("L99").Value = 1999.75
$ws.Range("M99").Value = 653
$ws.Range("N99").Value = -4995.75
$ws.Range("H107").Value = 1964.0834
$ws.Range("I107").Value = 1128.9412
$ws.Range("K107").Value = 1128.9412
$ws.Range("M107").Value = 791.0588
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 751.75
$ws.Range("I10").Value = 751.75
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 751.75
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -612.75
$ws.Range("H31").Value = 4806.846
$ws.Range("I31").Value = 2721
$ws.Range("J31").Value = 9500
$ws.Range("K31").Value = 2721
$ws.Range("L31").Value = 9500
$ws.Range("M31").Value = -2426
$ws.Range("N31").Value = -10090
$ws.Range("H34").Value = 4806.846
$ws.Range("I34").Value = 2721
$ws.Range("J34").Value = 9500
$ws.Range("K34").Value = 2721
$ws.Range("L34").Value = 9500
$ws.Range("M34").Value = -2519
$ws.Range("N34").Value = -9904
$ws.Range("H86").Value = 12734
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 18582
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 18582
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -20828
$ws.Range("H89").Value = 12734
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 18582
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 92910
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -104142
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1530.8379
$ws.Range("I5").Value = 1243.7084
$ws.Range("J5").Value = 2060.923
$ws.Range("K5").Value = 3731.1252
$ws.Range("L5").Value = 6182.768999999999
$ws.Range("M5").Value = -3619.1252
$ws.Range("N5").Value = -6406.768999999999
$ws.Range("H107").Value = 6491.7095
$ws.Range("I107").Value = 8912.227999999999
$ws.Range("J107").Value = 574.8889
$ws.Range("K107").Value = 26736.684
$ws.Range("L107").Value = 1724.6667
$ws.Range("M107").Value = -24816.684
$ws.Range("N107").Value = -5564.6667
$ws.Range("H122").Value = 1585.8182
$ws.Range("I122").Value = 496.66666
$ws.Range("J122").Value = 1757.7894
$ws.Range("K122").Value = 4469.99994
$ws.Range("L122").Value = 15820.1046
$ws.Range("M122").Value = -2019.99994
$ws.Range("N122").Value = -20720.1046
$ws.Range("H131").Value = 107182.945
$ws.Range("I131").Value = 555
$ws.Range("J131").Value = 114453.03
$ws.Range("K131").Value = 1665
$ws.Range("L131").Value = 343359.09
$ws.Range("M131").Value = 3375
$ws.Range("N131").Value = -353439.09
$ws.Range("H135").Value = 1530.8379
$ws.Range("I135").Value = 1243.7084
$ws.Range("J135").Value = 2060.923
$ws.Range("K135").Value = 11193.3756
$ws.Range("L135").Value = 18548.307
$ws.Range("M135").Value = -8658.375599999999
$ws.Range("N135").Value = -23618.307
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 903
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H132").Value = 6073800
$ws.Range("I132").Value = 8468327
$ws.Range("K132").Value = 25404981
$ws.Range("M132").Value = -25402451
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1834.5
$ws.Range("J82").Value = 2369
$ws.Range("L82").Value = 2369
$ws.Range("N82").Value = -3091
$ws.Range("H85").Value = 1834.5
$ws.Range("J85").Value = 2369
$ws.Range("L85").Value = 2369
$ws.Range("N85").Value = -4865
$ws.Range("H122").Value = 855748.1
$ws.Range("I122").Value = 1092194.9
$ws.Range("J122").Value = 4540
$ws.Range("K122").Value = 3276584.7
$ws.Range("L122").Value = 13620
$ws.Range("M122").Value = -3274134.7
$ws.Range("N122").Value = -18520
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 950
$ws.Range("I100").Value = 1200
$ws.Range("K100").Value = 2400
$ws.Range("M100").Value = -1859
$ws.Range("H109").Value = 27000
$ws.Range("J109").Value = 27000
$ws.Range("L109").Value = 27000
$ws.Range("N109").Value = -29774
$ws.Range("H136").Value = 29415194
$ws.Range("I136").Value = 45456384
$ws.Range("J136").Value = 6350.4165
$ws.Range("K136").Value = 136369152
$ws.Range("L136").Value = 19051.2495
$ws.Range("M136").Value = -136366602
$ws.Range("N136").Value = -24151.2495
